$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: copy the formatting of the last existing row (row 83) onto
# the new row 84, then set the new values on top of it.
$ws.Range("A83").Copy($ws.Range("A84"))

$ws.Range("A84").Value = 45884
$ws.Range("B84").Value = 0.2
